# Auto-generated edit script: apply IFRS figures fix for 흥아해운 workbook
# (commit message: "error solve ifrs list")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("D2").Value = 8251
$ws.Range("E2").Value = 186
$ws.Range("F2").Value = 186
$ws.Range("G2").Value = 209
$ws.Range("H2").Value = 183
$ws.Range("I2").Value = 183
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 6492
$ws.Range("L2").Value = 4818
$ws.Range("M2").Value = 1674
$ws.Range("N2").Value = 1674
$ws.Range("P2").Value = 424
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 259
$ws.Range("S2").Value = -367
$ws.Range("T2").Value = 78
$ws.Range("U2").Value = -73
$ws.Range("V2").Value = 2384
$ws.Range("W2").Value = 2.25
$ws.Range("X2").Value = 2.22
$ws.Range("Y2").Value = 11.37
$ws.Range("Z2").Value = 3.01
$ws.Range("AA2").Value = 287.83
$ws.Range("AB2").Value = 316.45
$ws.Range("AC2").Value = 396
$ws.Range("AD2").Value = 12.71
$ws.Range("AE2").Value = 3661
$ws.Range("AF2").Value = 1.38
$ws.Range("AG2").Value = 16
$ws.Range("AH2").Value = 0.33
$ws.Range("AI2").Value = 4.1
$ws.Range("AJ2").Value = 46292881
$ws.Range("O2").ClearContents()

# ---- Row 3 ----
$ws.Range("D3").Value = 8451
$ws.Range("E3").Value = 212
$ws.Range("F3").Value = 212
$ws.Range("G3").Value = 135
$ws.Range("H3").Value = 105
$ws.Range("I3").Value = 105
$ws.Range("K3").Value = 7916
$ws.Range("L3").Value = 6186
$ws.Range("M3").Value = 1731
$ws.Range("N3").Value = 1731
$ws.Range("P3").Value = 424
$ws.Range("Q3").Value = 339
$ws.Range("R3").Value = -6
$ws.Range("S3").Value = -295
$ws.Range("T3").Value = 278
$ws.Range("U3").Value = 62
$ws.Range("V3").Value = 2593
$ws.Range("W3").Value = 2.51
$ws.Range("X3").Value = 1.24
$ws.Range("Y3").Value = 6.18
$ws.Range("Z3").Value = 1.46
$ws.Range("AA3").Value = 357.39
$ws.Range("AB3").Value = 337.96
$ws.Range("AC3").Value = 227
$ws.Range("AD3").Value = 14.76
$ws.Range("AE3").Value = 3785
$ws.Range("AF3").Value = 0.89
$ws.Range("AG3").Value = 13
$ws.Range("AH3").Value = 0.38
$ws.Range("AI3").Value = 5.56
$ws.Range("AJ3").Value = 46292881
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# ---- Row 4 ----
$ws.Range("D4").Value = 8317
$ws.Range("E4").Value = 59
$ws.Range("F4").Value = 59
$ws.Range("G4").Value = -166
$ws.Range("H4").Value = -172
$ws.Range("I4").Value = -172
$ws.Range("K4").Value = 9585
$ws.Range("L4").Value = 7660
$ws.Range("M4").Value = 1925
$ws.Range("N4").Value = 1925
$ws.Range("P4").Value = 639
$ws.Range("Q4").Value = -14
$ws.Range("R4").Value = -147
$ws.Range("S4").Value = 103
$ws.Range("T4").Value = 330
$ws.Range("U4").Value = -344
$ws.Range("V4").Value = 2573
$ws.Range("W4").Value = 0.71
$ws.Range("X4").Value = -2.07
$ws.Range("Y4").Value = -9.4
$ws.Range("Z4").Value = -1.96
$ws.Range("AA4").Value = 397.97
$ws.Range("AB4").Value = 230.56
$ws.Range("AC4").Value = -324
$ws.Range("AD4").Value = -7.58
$ws.Range("AE4").Value = 2955
$ws.Range("AF4").Value = 0.83
$ws.Range("AG4").Value = 10
$ws.Range("AH4").Value = 0.4
$ws.Range("AI4").Value = -3.68
$ws.Range("AJ4").Value = 65701663
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# ---- Row 5 ----
$ws.Range("D5").Value = 8364
$ws.Range("E5").Value = -130
$ws.Range("F5").Value = -130
$ws.Range("G5").Value = -598
$ws.Range("H5").Value = -620
$ws.Range("I5").Value = -620
$ws.Range("K5").Value = 8734
$ws.Range("L5").Value = 7374
$ws.Range("M5").Value = 1360
$ws.Range("N5").Value = 1360
$ws.Range("P5").Value = 639
$ws.Range("Q5").Value = -69
$ws.Range("R5").Value = -20
$ws.Range("S5").Value = 145
$ws.Range("T5").Value = 465
$ws.Range("U5").Value = -534
$ws.Range("V5").Value = 2695
$ws.Range("W5").Value = -1.55
$ws.Range("X5").Value = -7.41
$ws.Range("Y5").Value = -37.75
$ws.Range("Z5").Value = -6.77
$ws.Range("AA5").Value = 542.29
$ws.Range("AB5").Value = 151.64
$ws.Range("AC5").Value = -944
$ws.Range("AD5").Value = -1.57
$ws.Range("AE5").Value = 2088
$ws.Range("AF5").Value = 0.71
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 65701663
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# ---- Row 6 ----
$ws.Range("D6").Value = 7539
$ws.Range("E6").Value = -376
$ws.Range("F6").Value = -376
$ws.Range("G6").Value = -912
$ws.Range("H6").Value = -865
$ws.Range("I6").Value = -865
$ws.Range("K6").Value = 8037
$ws.Range("L6").Value = 7306
$ws.Range("M6").Value = 731
$ws.Range("N6").Value = 731
$ws.Range("P6").Value = 878
$ws.Range("Q6").Value = -338
$ws.Range("R6").Value = 850
$ws.Range("S6").Value = -553
$ws.Range("T6").Value = 112
$ws.Range("U6").Value = -450
$ws.Range("V6").Value = 2540
$ws.Range("W6").Value = -4.99
$ws.Range("X6").Value = -11.48
$ws.Range("Y6").Value = -82.76000000000001
$ws.Range("Z6").Value = -10.32
$ws.Range("AA6").Value = 998.73
$ws.Range("AB6").Value = 13
$ws.Range("AC6").Value = -1097
$ws.Range("AD6").Value = -0.72
$ws.Range("AE6").Value = 894
$ws.Range("AF6").Value = 0.89
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 82376352

# ---- Rows 7-9: drop all figures except the first three columns (A-C) ----
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()

